$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "Project name" header text that was duplicated in D3,
# and fill in Project Name values for rows 3,4,6,7,8. Also mark row 3
# as Approved ("Yes").
$ws.Range("D3").Value = "Mass-Mass Stoichiometry Problem"
$ws.Range("G3").Value = "Yes"
$ws.Range("D4").Value = "Fourier Series"
$ws.Range("D6").Value = "Medical Imaging Applications"
$ws.Range("D7").Value = "Lattice Boltzmann Solvers"
$ws.Range("D8").Value = "Light"

# Update row heights for rows 3 and 6 to accommodate wrapped text.
$ws.Rows(3).RowHeight = 30
$ws.Rows(6).RowHeight = 30

# Update the active selection to reflect where the author left off editing.
$ws.Range("D9").Select()
